$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data feed re-scraped this round; several already-recorded
# fixtures shifted position in the sheet (new matches were inserted earlier
# in the scrape order), so the home/away/odds/url columns (F:V) for a handful
# of existing rows need to be re-paired, and two brand-new fixtures append
# at the bottom (rows 70-71). The Indice (A) and data_partida (E) columns for
# the existing rows are untouched - only F:V move.

# --- Re-pair match data (columns F:V) for rows whose fixture shifted ---
# Row 2 now holds the fixture previously recorded in row 3
$ws.Cells.Item(2, 6).Value = 'Stirling'
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 'Edinburgh City'
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 1.61
$ws.Cells.Item(2, 11).Value = '04/08/2023 16:33'
$ws.Cells.Item(2, 12).Value = 1.8
$ws.Cells.Item(2, 13).Value = '05/08/2023 15:59'
$ws.Cells.Item(2, 14).Value = 4.1
$ws.Cells.Item(2, 15).Value = '04/08/2023 16:33'
$ws.Cells.Item(2, 16).Value = 3.89
$ws.Cells.Item(2, 17).Value = '05/08/2023 15:59'
$ws.Cells.Item(2, 18).Value = 5.04
$ws.Cells.Item(2, 19).Value = '04/08/2023 16:33'
$ws.Cells.Item(2, 20).Value = 4.08
$ws.Cells.Item(2, 21).Value = '05/08/2023 15:59'
$ws.Cells.Item(2, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/stirling-edinburgh-city/dE0H0tig/'

# Row 3 now holds the fixture previously recorded in row 2
$ws.Cells.Item(3, 6).Value = 'Queen of South'
$ws.Cells.Item(3, 7).Value = 3
$ws.Cells.Item(3, 8).Value = 'Alloa'
$ws.Cells.Item(3, 9).Value = 4
$ws.Cells.Item(3, 10).Value = 1.84
$ws.Cells.Item(3, 11).Value = '03/08/2023 09:13'
$ws.Cells.Item(3, 12).Value = 1.96
$ws.Cells.Item(3, 13).Value = '05/08/2023 15:48'
$ws.Cells.Item(3, 14).Value = 3.55
$ws.Cells.Item(3, 15).Value = '03/08/2023 09:13'
$ws.Cells.Item(3, 16).Value = 3.61
$ws.Cells.Item(3, 17).Value = '05/08/2023 15:48'
$ws.Cells.Item(3, 18).Value = 3.75
$ws.Cells.Item(3, 19).Value = '03/08/2023 09:13'
$ws.Cells.Item(3, 20).Value = 3.68
$ws.Cells.Item(3, 21).Value = '05/08/2023 15:48'
$ws.Cells.Item(3, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/queen-of-south-alloa/IN1D1Mym/'

# Row 42 now holds the fixture previously recorded in row 43
$ws.Cells.Item(42, 6).Value = 'Montrose'
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(42, 8).Value = 'Alloa'
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = 1.96
$ws.Cells.Item(42, 11).Value = '05/10/2023 08:12'
$ws.Cells.Item(42, 12).Value = 2.02
$ws.Cells.Item(42, 13).Value = '07/10/2023 15:03'
$ws.Cells.Item(42, 14).Value = 3.48
$ws.Cells.Item(42, 15).Value = '05/10/2023 08:12'
$ws.Cells.Item(42, 16).Value = 3.63
$ws.Cells.Item(42, 17).Value = '07/10/2023 15:03'
$ws.Cells.Item(42, 18).Value = 3.28
$ws.Cells.Item(42, 19).Value = '05/10/2023 08:12'
$ws.Cells.Item(42, 20).Value = 3.48
$ws.Cells.Item(42, 21).Value = '07/10/2023 15:03'
$ws.Cells.Item(42, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/montrose-alloa/CC1oI6Xm/'

# Row 43 now holds the fixture previously recorded in row 42
$ws.Cells.Item(43, 6).Value = 'Queen of South'
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 'Stirling'
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 1.96
$ws.Cells.Item(43, 11).Value = '05/10/2023 08:12'
$ws.Cells.Item(43, 12).Value = 2.11
$ws.Cells.Item(43, 13).Value = '07/10/2023 15:41'
$ws.Cells.Item(43, 14).Value = 3.49
$ws.Cells.Item(43, 15).Value = '05/10/2023 08:12'
$ws.Cells.Item(43, 16).Value = 3.6
$ws.Cells.Item(43, 17).Value = '07/10/2023 15:41'
$ws.Cells.Item(43, 18).Value = 3.27
$ws.Cells.Item(43, 19).Value = '05/10/2023 08:12'
$ws.Cells.Item(43, 20).Value = 3.26
$ws.Cells.Item(43, 21).Value = '07/10/2023 15:41'
$ws.Cells.Item(43, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/queen-of-south-stirling/QF5kHQnf/'

# Row 51 now holds the fixture previously recorded in row 52
$ws.Cells.Item(51, 6).Value = 'Kelty Hearts'
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 'Stirling'
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 2.07
$ws.Cells.Item(51, 11).Value = '26/10/2023 09:12'
$ws.Cells.Item(51, 12).Value = 2.44
$ws.Cells.Item(51, 13).Value = '28/10/2023 15:53'
$ws.Cells.Item(51, 14).Value = 3.41
$ws.Cells.Item(51, 15).Value = '26/10/2023 09:12'
$ws.Cells.Item(51, 16).Value = 3.62
$ws.Cells.Item(51, 17).Value = '28/10/2023 15:53'
$ws.Cells.Item(51, 18).Value = 3.08
$ws.Cells.Item(51, 19).Value = '26/10/2023 09:12'
$ws.Cells.Item(51, 20).Value = 2.69
$ws.Cells.Item(51, 21).Value = '28/10/2023 15:53'
$ws.Cells.Item(51, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/kelty-hearts-stirling/pQLOVP1m/'

# Row 52 now holds the fixture previously recorded in row 51
$ws.Cells.Item(52, 6).Value = 'Queen of South'
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 'Cove Rangers'
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 2.37
$ws.Cells.Item(52, 11).Value = '27/10/2023 04:13'
$ws.Cells.Item(52, 12).Value = 2.39
$ws.Cells.Item(52, 13).Value = '28/10/2023 15:43'
$ws.Cells.Item(52, 14).Value = 3.5
$ws.Cells.Item(52, 15).Value = '27/10/2023 04:13'
$ws.Cells.Item(52, 16).Value = 3.83
$ws.Cells.Item(52, 17).Value = '28/10/2023 15:34'
$ws.Cells.Item(52, 18).Value = 2.6
$ws.Cells.Item(52, 19).Value = '27/10/2023 04:13'
$ws.Cells.Item(52, 20).Value = 2.65
$ws.Cells.Item(52, 21).Value = '28/10/2023 15:43'
$ws.Cells.Item(52, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/queen-of-south-cove-rangers/thFXT3Wa/'

# Row 53 now holds the fixture previously recorded in row 54
$ws.Cells.Item(53, 6).Value = 'Falkirk'
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 'Alloa'
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 1.28
$ws.Cells.Item(53, 11).Value = '26/10/2023 09:12'
$ws.Cells.Item(53, 12).Value = 1.28
$ws.Cells.Item(53, 13).Value = '28/10/2023 15:33'
$ws.Cells.Item(53, 14).Value = 5.23
$ws.Cells.Item(53, 15).Value = '26/10/2023 09:12'
$ws.Cells.Item(53, 16).Value = 5.68
$ws.Cells.Item(53, 17).Value = '28/10/2023 15:59'
$ws.Cells.Item(53, 18).Value = 7.7
$ws.Cells.Item(53, 19).Value = '26/10/2023 09:12'
$ws.Cells.Item(53, 20).Value = 10.24
$ws.Cells.Item(53, 21).Value = '28/10/2023 15:59'
$ws.Cells.Item(53, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/falkirk-alloa/6ZMKW5ns/'

# Row 54 now holds the fixture previously recorded in row 53
$ws.Cells.Item(54, 6).Value = 'Edinburgh City'
$ws.Cells.Item(54, 7).Value = 3
$ws.Cells.Item(54, 8).Value = 'Annan'
$ws.Cells.Item(54, 9).Value = 2
$ws.Cells.Item(54, 10).Value = 2.94
$ws.Cells.Item(54, 11).Value = '26/10/2023 09:12'
$ws.Cells.Item(54, 12).Value = 3.12
$ws.Cells.Item(54, 13).Value = '28/10/2023 15:41'
$ws.Cells.Item(54, 14).Value = 3.48
$ws.Cells.Item(54, 15).Value = '26/10/2023 09:12'
$ws.Cells.Item(54, 16).Value = 4.09
$ws.Cells.Item(54, 17).Value = '28/10/2023 15:59'
$ws.Cells.Item(54, 18).Value = 2.11
$ws.Cells.Item(54, 19).Value = '26/10/2023 09:12'
$ws.Cells.Item(54, 20).Value = 2.03
$ws.Cells.Item(54, 21).Value = '28/10/2023 15:41'
$ws.Cells.Item(54, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/edinburgh-city-annan/lfzgynvQ/'

# Row 61 now holds the fixture previously recorded in row 65
$ws.Cells.Item(61, 6).Value = 'Montrose'
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 'Annan'
$ws.Cells.Item(61, 9).Value = 1
$ws.Cells.Item(61, 10).Value = 1.9
$ws.Cells.Item(61, 11).Value = '09/11/2023 09:12'
$ws.Cells.Item(61, 12).Value = 1.79
$ws.Cells.Item(61, 13).Value = '11/11/2023 15:54'
$ws.Cells.Item(61, 14).Value = 3.92
$ws.Cells.Item(61, 15).Value = '09/11/2023 09:12'
$ws.Cells.Item(61, 16).Value = 4.32
$ws.Cells.Item(61, 17).Value = '11/11/2023 15:54'
$ws.Cells.Item(61, 18).Value = 3.22
$ws.Cells.Item(61, 19).Value = '09/11/2023 09:12'
$ws.Cells.Item(61, 20).Value = 3.74
$ws.Cells.Item(61, 21).Value = '11/11/2023 15:54'
$ws.Cells.Item(61, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/montrose-annan/lzUFTfxH/'

# Row 62 now holds the fixture previously recorded in row 64
$ws.Cells.Item(62, 6).Value = 'Kelty Hearts'
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 'Cove Rangers'
$ws.Cells.Item(62, 9).Value = 1
$ws.Cells.Item(62, 10).Value = 2.46
$ws.Cells.Item(62, 11).Value = '09/11/2023 09:12'
$ws.Cells.Item(62, 12).Value = 3.06
$ws.Cells.Item(62, 13).Value = '11/11/2023 15:57'
$ws.Cells.Item(62, 14).Value = 3.33
$ws.Cells.Item(62, 15).Value = '09/11/2023 09:12'
$ws.Cells.Item(62, 16).Value = 3.51
$ws.Cells.Item(62, 17).Value = '11/11/2023 15:57'
$ws.Cells.Item(62, 18).Value = 2.53
$ws.Cells.Item(62, 19).Value = '09/11/2023 09:12'
$ws.Cells.Item(62, 20).Value = 2.24
$ws.Cells.Item(62, 21).Value = '11/11/2023 15:53'
$ws.Cells.Item(62, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/kelty-hearts-cove-rangers/AqVBUzNA/'

# Row 63 now holds the fixture previously recorded in row 62
$ws.Cells.Item(63, 6).Value = 'Falkirk'
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = 'Edinburgh City'
$ws.Cells.Item(63, 9).Value = 1
$ws.Cells.Item(63, 10).Value = 1.14
$ws.Cells.Item(63, 11).Value = '09/11/2023 09:12'
$ws.Cells.Item(63, 12).Value = 1.09
$ws.Cells.Item(63, 13).Value = '10/11/2023 16:03'
$ws.Cells.Item(63, 14).Value = 7.57
$ws.Cells.Item(63, 15).Value = '09/11/2023 09:12'
$ws.Cells.Item(63, 16).Value = 11.04
$ws.Cells.Item(63, 17).Value = '11/11/2023 15:46'
$ws.Cells.Item(63, 18).Value = 11.46
$ws.Cells.Item(63, 19).Value = '09/11/2023 09:12'
$ws.Cells.Item(63, 20).Value = 22.57
$ws.Cells.Item(63, 21).Value = '11/11/2023 15:46'
$ws.Cells.Item(63, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/falkirk-edinburgh-city/2JR3Wdhb/'

# Row 64 now holds the fixture previously recorded in row 61
$ws.Cells.Item(64, 6).Value = 'Alloa'
$ws.Cells.Item(64, 7).Value = 1
$ws.Cells.Item(64, 8).Value = 'Queen of South'
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2.41
$ws.Cells.Item(64, 11).Value = '09/11/2023 09:12'
$ws.Cells.Item(64, 12).Value = 2.46
$ws.Cells.Item(64, 13).Value = '11/11/2023 15:49'
$ws.Cells.Item(64, 14).Value = 3.29
$ws.Cells.Item(64, 15).Value = '09/11/2023 09:12'
$ws.Cells.Item(64, 16).Value = 3.52
$ws.Cells.Item(64, 17).Value = '11/11/2023 15:20'
$ws.Cells.Item(64, 18).Value = 2.62
$ws.Cells.Item(64, 19).Value = '09/11/2023 09:12'
$ws.Cells.Item(64, 20).Value = 2.72
$ws.Cells.Item(64, 21).Value = '11/11/2023 15:49'
$ws.Cells.Item(64, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/alloa-queen-of-south/KfxaXxwh/'

# Row 65 now holds the fixture previously recorded in row 63
$ws.Cells.Item(65, 6).Value = 'Hamilton'
$ws.Cells.Item(65, 7).Value = 5
$ws.Cells.Item(65, 8).Value = 'Stirling'
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 1.37
$ws.Cells.Item(65, 11).Value = '09/11/2023 09:12'
$ws.Cells.Item(65, 12).Value = 1.37
$ws.Cells.Item(65, 13).Value = '11/11/2023 15:52'
$ws.Cells.Item(65, 14).Value = 4.59
$ws.Cells.Item(65, 15).Value = '09/11/2023 09:12'
$ws.Cells.Item(65, 16).Value = 4.79
$ws.Cells.Item(65, 17).Value = '11/11/2023 15:52'
$ws.Cells.Item(65, 18).Value = 6.39
$ws.Cells.Item(65, 19).Value = '09/11/2023 09:12'
$ws.Cells.Item(65, 20).Value = 8.41
$ws.Cells.Item(65, 21).Value = '11/11/2023 15:52'
$ws.Cells.Item(65, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/hamilton-stirling/GMV7VG74/'

# Row 67 now holds the fixture previously recorded in row 68
$ws.Cells.Item(67, 6).Value = 'Cove Rangers'
$ws.Cells.Item(67, 7).Value = 3
$ws.Cells.Item(67, 8).Value = 'Stirling'
$ws.Cells.Item(67, 9).Value = 1
$ws.Cells.Item(67, 10).Value = 1.68
$ws.Cells.Item(67, 11).Value = '16/11/2023 09:13'
$ws.Cells.Item(67, 12).Value = 1.66
$ws.Cells.Item(67, 13).Value = '18/11/2023 15:54'
$ws.Cells.Item(67, 14).Value = 3.77
$ws.Cells.Item(67, 15).Value = '16/11/2023 09:13'
$ws.Cells.Item(67, 16).Value = 3.75
$ws.Cells.Item(67, 17).Value = '18/11/2023 15:54'
$ws.Cells.Item(67, 18).Value = 4.15
$ws.Cells.Item(67, 19).Value = '16/11/2023 09:13'
$ws.Cells.Item(67, 20).Value = 5.31
$ws.Cells.Item(67, 21).Value = '18/11/2023 15:54'
$ws.Cells.Item(67, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/cove-rangers-stirling/8xpHRY6T/'

# Row 68 now holds the fixture previously recorded in row 67
$ws.Cells.Item(68, 6).Value = 'Edinburgh City'
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 'Alloa'
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 3.04
$ws.Cells.Item(68, 11).Value = '16/11/2023 09:13'
$ws.Cells.Item(68, 12).Value = 4.08
$ws.Cells.Item(68, 13).Value = '18/11/2023 15:15'
$ws.Cells.Item(68, 14).Value = 3.52
$ws.Cells.Item(68, 15).Value = '16/11/2023 09:13'
$ws.Cells.Item(68, 16).Value = 3.97
$ws.Cells.Item(68, 17).Value = '18/11/2023 15:15'
$ws.Cells.Item(68, 18).Value = 2.05
$ws.Cells.Item(68, 19).Value = '16/11/2023 09:13'
$ws.Cells.Item(68, 20).Value = 1.78
$ws.Cells.Item(68, 21).Value = '18/11/2023 15:15'
$ws.Cells.Item(68, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/edinburgh-city-alloa/QHWrNhFp/'

# --- Append the two newly scraped fixtures as rows 70 and 71 ---
$ws.Range("A69:V69").Copy()
$ws.Range("A70").PasteSpecial(-4122)
$ws.Range("A71").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 70
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = 'scotland'
$ws.Cells.Item(70, 3).Value = 'league-one'
$ws.Cells.Item(70, 4).Value = '2023-2024'
$ws.Cells.Item(70, 5).Value = 45258.86458333334
$ws.Cells.Item(70, 6).Value = 'Annan'
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 'Hamilton'
$ws.Cells.Item(70, 9).Value = 2
$ws.Cells.Item(70, 10).Value = 5.04
$ws.Cells.Item(70, 11).Value = '24/11/2023 01:13'
$ws.Cells.Item(70, 12).Value = 5.51
$ws.Cells.Item(70, 13).Value = '28/11/2023 20:44'
$ws.Cells.Item(70, 14).Value = 4.22
$ws.Cells.Item(70, 15).Value = '24/11/2023 01:13'
$ws.Cells.Item(70, 16).Value = 4.31
$ws.Cells.Item(70, 17).Value = '28/11/2023 20:44'
$ws.Cells.Item(70, 18).Value = 1.5
$ws.Cells.Item(70, 19).Value = '24/11/2023 01:13'
$ws.Cells.Item(70, 20).Value = 1.55
$ws.Cells.Item(70, 21).Value = '28/11/2023 20:36'
$ws.Cells.Item(70, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/annan-hamilton/UooDSEiN/'

# Row 71
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = 'scotland'
$ws.Cells.Item(71, 3).Value = 'league-one'
$ws.Cells.Item(71, 4).Value = '2023-2024'
$ws.Cells.Item(71, 5).Value = 45258.86458333334
$ws.Cells.Item(71, 6).Value = 'Montrose'
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 'Falkirk'
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 5.96
$ws.Cells.Item(71, 11).Value = '24/11/2023 01:13'
$ws.Cells.Item(71, 12).Value = 8.48
$ws.Cells.Item(71, 13).Value = '28/11/2023 19:24'
$ws.Cells.Item(71, 14).Value = 4.38
$ws.Cells.Item(71, 15).Value = '24/11/2023 01:13'
$ws.Cells.Item(71, 16).Value = 5.02
$ws.Cells.Item(71, 17).Value = '28/11/2023 19:24'
$ws.Cells.Item(71, 18).Value = 1.44
$ws.Cells.Item(71, 19).Value = '24/11/2023 01:13'
$ws.Cells.Item(71, 20).Value = 1.35
$ws.Cells.Item(71, 21).Value = '28/11/2023 19:24'
$ws.Cells.Item(71, 22).Value = 'https://www.betexplorer.com/football/scotland/league-one/montrose-falkirk/vLzmMCUj/'
